$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ G = 1.609405; H = 4.828215; I = 0.1374279556489526; J = 0.1582737093407849; M = 0.690981; O = 0.9906161271292202; P = 0.993724455114346; Q = 1.112068276305; R = 10.008614486745; S = 0.1361383491842516; T = 0.1572804555735979 }
    3  = @{ G = 1.609405; H = 4.828215; I = 0.1374279556489526; J = 0.1582737093407849; K = 1; L = 0.5; M = 0.0065455; N = 0.013091; O = 0.009383872870779821; P = 0.006275544885653829; Q = 0.0105343604275; R = 0.063206162565; S = 0.001289606464700939; T = 0.0009932537671870237 }
    4  = @{ I = 0.4033933174334258; J = 0.4645820158786568; M = 0.690981; O = 0.9906161271292202; P = 0.993724455114346; Q = 3.264262420791999; S = 0.3996079258257084; T = 0.4616665105849427 }
    5  = @{ I = 0.4033933174334258; J = 0.4645820158786568; K = 1; L = 0.5; M = 0.0065455; N = 0.013091; O = 0.009383872870779821; P = 0.006275544885653829; Q = 0.03092158782266666; R = 0.185529526936; S = 0.003785391607717397; T = 0.002915505293714051 }
    6  = @{ G = 0.476314; H = 1.428942; I = 0.04067270778143176; J = 0.04684214575631779; M = 0.690981; O = 0.9906161271292202; P = 0.993724455114346; Q = 0.3291239240339999; R = 2.962115316306; S = 0.04029104026230043; T = 0.04654818576808367 }
    7  = @{ G = 0.476314; H = 1.428942; I = 0.04067270778143176; J = 0.04684214575631779; K = 1; L = 0.5; M = 0.0065455; N = 0.013091; O = 0.009383872870779821; P = 0.006275544885653829; Q = 0.003117713287; R = 0.018706279722; S = 0.0003816675191313328; T = 0.0002939599882341113 }
    8  = @{ G = 4.627222; H = 9.254443999999999; I = 0.3951209669373822; J = 0.3033699161629238; M = 0.690981; O = 0.9906161271292202; P = 0.993724455114346; Q = 3.197322484782; R = 19.183934908692; S = 0.3914132020150622; T = 0.3014661046370863 }
    9  = @{ G = 4.627222; H = 9.254443999999999; I = 0.3951209669373822; J = 0.3033699161629238; K = 1; L = 0.5; M = 0.0065455; N = 0.013091; O = 0.009383872870779821; P = 0.006275544885653829; Q = 0.030287481601; R = 0.121149926404; S = 0.003707764922319991; T = 0.001903811525837467 }
    10 = @{ E = 2; F = 0.6666666666666666; G = 0.27386; H = 0.82158; I = 0.02338505219880773; J = 0.02693221286131667; M = 0.690981; O = 0.9906161271292202; P = 0.993724455114346; Q = 0.18923205666; R = 1.70308850994; S = 0.02316560984189757; T = 0.02676319855063549 }
    11 = @{ E = 2; F = 0.6666666666666666; G = 0.27386; H = 0.82158; I = 0.02338505219880773; J = 0.02693221286131667; K = 1; L = 0.5; M = 0.0065455; N = 0.013091; O = 0.009383872870779821; P = 0.006275544885653829; Q = 0.00179255063; R = 0.01075530378; S = 0.0002194423569101618; T = 0.0001690143106811761 }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
